# Weekly refresh of the "Hortaliza, Vega Monumental Concepción - Pepino dulce" data.
# The underlying table (rows 2-39) is re-sorted by date: each row's
# D (Fecha), I (Calidad), J (Volumen), K (Precio mínimo), L (Precio máximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) values are replaced by
# those of another row in the same table (a pure row permutation - the
# other columns, which are constant for this sub-sheet, are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 39

# For output row (index into array, 0-based corresponds to $firstRow), the
# value is taken from this source row of the *original* sheet.
$srcRow = @(28, 29, 19, 20, 36, 37, 23, 24, 6, 7, 10, 11, 38, 39, 25, 26, 16, 17, 32, 33, 34, 35, 14, 15, 21, 22, 2, 3, 18, 8, 9, 12, 13, 30, 31, 27, 4, 5)

$cols = @(4, 9, 10, 11, 12, 13, 16)   # D, I, J, K, L, M, P

# Snapshot the original values for every touched cell before writing
# anything, since several source/destination rows overlap.
$orig = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $orig["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

for ($i = 0; $i -lt $srcRow.Length; $i++) {
    $destRow = $firstRow + $i
    $source = $srcRow[$i]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $orig["$source,$c"]
    }
}
